$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the course id casing for John Smith's course: Data200 -> DATA201 -> DATA200
# (mirrors the natural shared-string table growth of the original edit)
$ws.Range("D2").Value = "DATA201"
$ws.Range("D2").Value = "DATA200"

# Add a new professor row for Ken Douglas (Visiting Lecturer, DATA201), filled right to left
$ws.Range("D3").Value = "DATA201"
$ws.Range("C3").Value = "Visiting Lecturer"
$ws.Range("B3").Value = "Ken Douglas"
$ws.Range("A3").Value = "douglas@myschool.edu"

# Hyperlink the new professor's email, mirroring A2's mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:douglas@myschool.edu")

# Restore the Hyperlink cell style so A3 matches A2's formatting exactly
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B5").Select()
